$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original sheet, add the new "rotation" sheet after it
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Distance"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "rotation"

# ---------------------------------------------------------------------
# 2. Populate "rotation" with a copy of the Run/Power/Rotations data
#    (same A/B/C layout that "Distance" already has)
# ---------------------------------------------------------------------
$ws1.Range("A1:C31").Copy()
$ws2.Range("A1").PasteSpecial()

# Header row text (new labels "Rotations" / "Total Angle", and the
# "In/Pwr" helper-column labels shifted one column to the right)
$ws2.Range("C1").Value = "Rotations"
$ws2.Range("D1").Value = "Total Angle"
$ws2.Range("F1").Value = "In/Pwr"
$ws2.Range("G1").Value = "In/Pwr (Corrected)"

# "Total Angle" = 360 * Rotations, computed only on the rows that have
# a Rotations entry (same sparsely-filled pattern as column C)
$ws2.Range("D2").Formula = "=360*C2"
$ws2.Range("D7").Formula = "=360*C7"
$ws2.Range("D12").Formula = "=360*C12"
$ws2.Range("D17").Formula = "=360*C17"
$ws2.Range("D22").Formula = "=360*C22"
$ws2.Range("D27").Formula = "=360*C27"

# Helper slope formulas (mirrors Distance!E2/F2, shifted right one
# column and pointing at the new Total Angle column)
$ws2.Range("F2").Formula = "=SLOPE(D:D,B:B)"
$ws2.Range("G2").Formula = "=SLOPE(D7:D31,B7:B31)"
$ws2.Range("G2").NumberFormat = $ws1.Range("F2").NumberFormat
$ws2.Range("H2").NumberFormat = $ws1.Range("G2").NumberFormat

$wb.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. New "Angle vs. Power" scatter chart on the "rotation" sheet
# ---------------------------------------------------------------------
$co = $ws2.ChartObjects().Add(323850, 700000, 4500000, 2800000)
$co.Name = "Chart 1"
$chart = $co.Chart
$chart.ChartType = -4169
$chart.SetSourceData($ws2.Range("B1:B31"))

$series = $chart.SeriesCollection(1)
$series.Name = "=rotation!`$D`$1"
$series.XValues = $ws2.Range("B2:B31")
$series.Values = $ws2.Range("D2:D31")
$series.Trendlines().Add()

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Angle vs. Power"

Write-Output "edit applied"
